$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices) are written as literal text,
# matching the source data which stores them as strings (e.g. "10.80", "70.677.77").
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.683.62"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.519.62"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "619.45"
$ws.Range("E5").Value = "  +2.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.98"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.515.54"
$ws.Range("E8").Value = "  -2.22%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  -2.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.06"
$ws.Range("E11").Value = "  -5.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.587"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.46"
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.084.18"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.38"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "608.71"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.519.11"
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.736.92"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.122"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.71"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.882"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.18"
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.64"
$ws.Range("E24").Value = "  -3.26%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.82"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("E26").Value = "  -2.44%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -4.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.69"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.06"
$ws.Range("E30").Value = "  -3.07%  "
$ws.Range("E31").Value = "  -3.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.09"
$ws.Range("E32").Value = "  -5.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.29"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "639.88"
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.82"
$ws.Range("E35").Value = "  -6.16%  "
$ws.Range("E36").Value = "  -2.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.80"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0485"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.43"
$ws.Range("E39").Value = "  -8.80%  "
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "56.60"
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.350.66"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0717"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.312"
$ws.Range("E45").Value = "  -4.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.95"
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "31.78"
$ws.Range("E47").Value = "  -4.08%  "
$ws.Range("E48").Value = "  -6.65%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.25"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("E51").Value = "  -0.04%  "
